# #141 update changelog, correct TokenData.csv
#
# Appends a new changelog entry (row 22) to the "Tabelle1" sheet:
#   A22 = 2021-08-07 (serial 44415)
#   B22 = "1.5.3"
#   C22 = multi-line release notes
# mirroring the date / version / wrapped-text formatting already used by
# the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate    = 44415
$newVersion = "1.5.3"
$newChanges = "Update graphics for new USDC-pool. The following evaluations were changed`n- Overview`n- Coins`n- Coinprices`n- Volume`n- Price stability`n- TVL`n- Liquidity Token`n- Fees`n- Cryptos-DAT"

# Copy the previous row's cell formats (date format / version style /
# wrapped-text style) onto the new row before writing values, so the new
# cells reuse the workbook's existing styles instead of minting new ones.
$ws.Range("A21:C21").Copy() | Out-Null
$ws.Range("A22:C22").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("A22").Value = $newDate
$ws.Range("B22").Value = $newVersion
$ws.Range("C22").Value = $newChanges

# The long, wrapped text needs a taller row, same as the other
# multi-line entries above it.
$ws.Rows.Item(22).RowHeight = 150

# Leave the selection where Excel would after typing the new row.
$ws.Range("C23").Select() | Out-Null
